$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "64.259.77"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.56%  "
Set-TextValue $ws.Cells.Item(3, 4) "3.504.21"
Set-TextValue $ws.Cells.Item(3, 5) "  -0.29%  "
Set-TextValue $ws.Cells.Item(4, 5) "  -0.01%  "
Set-TextValue $ws.Cells.Item(5, 4) "584.61"
Set-TextValue $ws.Cells.Item(5, 5) "  -0.39%  "
Set-TextValue $ws.Cells.Item(6, 4) "135.11"
Set-TextValue $ws.Cells.Item(6, 5) "  +1.60%  "
Set-TextValue $ws.Cells.Item(7, 4) "3.504.51"
Set-TextValue $ws.Cells.Item(7, 5) "  -0.28%  "
Set-TextValue $ws.Cells.Item(8, 5) "  -0.01%  "
Set-TextValue $ws.Cells.Item(9, 5) "  -0.50%  "
Set-TextValue $ws.Cells.Item(10, 5) "  +0.02%  "
Set-TextValue $ws.Cells.Item(11, 4) "7.12"
Set-TextValue $ws.Cells.Item(11, 5) "  -0.50%  "
Set-TextValue $ws.Cells.Item(12, 5) "  -3.81%  "
Set-TextValue $ws.Cells.Item(13, 4) "4.099.22"
Set-TextValue $ws.Cells.Item(13, 5) "  -0.35%  "
Set-TextValue $ws.Cells.Item(14, 5) "  -0.79%  "
Set-TextValue $ws.Cells.Item(15, 5) "  +1.11%  "
Set-TextValue $ws.Cells.Item(16, 4) "3.504.88"
Set-TextValue $ws.Cells.Item(16, 5) "  -0.23%  "
Set-TextValue $ws.Cells.Item(17, 4) "26.38"
Set-TextValue $ws.Cells.Item(17, 5) "  -5.57%  "
Set-TextValue $ws.Cells.Item(18, 4) "64.269.62"
Set-TextValue $ws.Cells.Item(18, 5) "  -0.54%  "
Set-TextValue $ws.Cells.Item(19, 4) "9.76"
Set-TextValue $ws.Cells.Item(19, 5) "  -2.48%  "
Set-TextValue $ws.Cells.Item(20, 4) "13.88"
Set-TextValue $ws.Cells.Item(20, 5) "  -2.83%  "
Set-TextValue $ws.Cells.Item(21, 4) "5.59"
Set-TextValue $ws.Cells.Item(21, 5) "  -2.19%  "
Set-TextValue $ws.Cells.Item(22, 4) "383.77"
Set-TextValue $ws.Cells.Item(22, 5) "  -2.45%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.570"
Set-TextValue $ws.Cells.Item(23, 5) "  -1.65%  "
Set-TextValue $ws.Cells.Item(24, 4) "3.640.33"
Set-TextValue $ws.Cells.Item(24, 5) "  -0.41%  "
Set-TextValue $ws.Cells.Item(25, 4) "73.97"
Set-TextValue $ws.Cells.Item(25, 5) "  -0.46%  "
Set-TextValue $ws.Cells.Item(26, 5) "  +0.00%  "
Set-TextValue $ws.Cells.Item(28, 5) "  +3.30%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.58"
Set-TextValue $ws.Cells.Item(29, 5) "  -0.30%  "
Set-TextValue $ws.Cells.Item(30, 4) "7.56"
Set-TextValue $ws.Cells.Item(30, 5) "  +0.39%  "
Set-TextValue $ws.Cells.Item(31, 5) "  +0.00%  "
Set-TextValue $ws.Cells.Item(32, 4) "8.30"
Set-TextValue $ws.Cells.Item(32, 5) "  +0.60%  "
Set-TextValue $ws.Cells.Item(33, 5) "  -1.73%  "
Set-TextValue $ws.Cells.Item(34, 4) "3.522.35"
Set-TextValue $ws.Cells.Item(35, 5) "  -0.02%  "
Set-TextValue $ws.Cells.Item(36, 2) "Kaspa"
Set-TextValue $ws.Cells.Item(36, 3) "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Cells.Item(36, 4) "0.146"
Set-TextValue $ws.Cells.Item(36, 5) "  -0.24%  "
Set-TextValue $ws.Cells.Item(37, 2) "EthereumClassic"
Set-TextValue $ws.Cells.Item(37, 3) "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Cells.Item(37, 4) "23.57"
Set-TextValue $ws.Cells.Item(37, 5) "  -2.00%  "
Set-TextValue $ws.Cells.Item(38, 4) "5.34"
Set-TextValue $ws.Cells.Item(38, 5) "  +1.00%  "
Set-TextValue $ws.Cells.Item(39, 4) "6.86"
Set-TextValue $ws.Cells.Item(39, 5) "  -1.81%  "
Set-TextValue $ws.Cells.Item(40, 5) "  -4.72%  "
Set-TextValue $ws.Cells.Item(41, 4) "164.71"
Set-TextValue $ws.Cells.Item(41, 5) "  -3.84%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.0784"
Set-TextValue $ws.Cells.Item(42, 5) "  -3.62%  "
Set-TextValue $ws.Cells.Item(43, 5) "  -0.67%  "
Set-TextValue $ws.Cells.Item(44, 4) "26.17"
Set-TextValue $ws.Cells.Item(44, 5) "  -2.00%  "
Set-TextValue $ws.Cells.Item(45, 5) "  +0.01%  "
Set-TextValue $ws.Cells.Item(46, 5) "  -1.13%  "
Set-TextValue $ws.Cells.Item(47, 5) "  -0.10%  "
Set-TextValue $ws.Cells.Item(48, 5) "  -0.85%  "
Set-TextValue $ws.Cells.Item(49, 5) "  -2.25%  "
Set-TextValue $ws.Cells.Item(50, 4) "2.477.10"
Set-TextValue $ws.Cells.Item(50, 5) "  +0.06%  "
Set-TextValue $ws.Cells.Item(51, 5) "  +1.05%  "
